# Auto-generated script applying numeric corrections to leve profit calculations
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 300
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -188
$ws.Range("N6").Value = -1724

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2111.92
$ws.Range("I40").Value = 1824.3572
$ws.Range("J40").Value = 2477.9092
$ws.Range("K40").Value = 1824.3572
$ws.Range("L40").Value = 2477.9092
$ws.Range("M40").Value = -1649.3572
$ws.Range("N40").Value = -2827.9092

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 37039896
$ws.Range("I64").Value = 52633852
$ws.Range("J64").Value = 4246.25
$ws.Range("K64").Value = 52633852
$ws.Range("L64").Value = 4246.25
$ws.Range("M64").Value = -52633604
$ws.Range("N64").Value = -4742.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 37039896
$ws.Range("I67").Value = 52633852
$ws.Range("J67").Value = 4246.25
$ws.Range("K67").Value = 52633852
$ws.Range("L67").Value = 4246.25
$ws.Range("M67").Value = -52632994
$ws.Range("N67").Value = -5962.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1800
$ws.Range("I94").Value = 1700
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1700
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -1249
$ws.Range("N94").Value = -2902

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1203.5555
$ws.Range("I101").Value = 805.3333
$ws.Range("J101").Value = 2000
$ws.Range("K101").Value = 2415.9999
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = -793.9998999999998
$ws.Range("N101").Value = -9244

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 19317.363
$ws.Range("I116").Value = 5834.1665
$ws.Range("J116").Value = 35497.2
$ws.Range("K116").Value = 5834.1665
$ws.Range("L116").Value = 35497.2
$ws.Range("M116").Value = -2392.1665
$ws.Range("N116").Value = -42381.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 756.21277
$ws.Range("I2").Value = 596.439
$ws.Range("J2").Value = 1848
$ws.Range("K2").Value = 596.439
$ws.Range("L2").Value = 1848
$ws.Range("M2").Value = -483.439
$ws.Range("N2").Value = -2074

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 14990
$ws.Range("I55").Value = 14980
$ws.Range("J55").Value = 15000
$ws.Range("K55").Value = 14980
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -14665
$ws.Range("N55").Value = -15630

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3432.7693
$ws.Range("I88").Value = 2953.5
$ws.Range("J88").Value = 4199.6
$ws.Range("K88").Value = 2953.5
$ws.Range("L88").Value = 4199.6
$ws.Range("M88").Value = -2547.5
$ws.Range("N88").Value = -5011.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3432.7693
$ws.Range("I91").Value = 2953.5
$ws.Range("J91").Value = 4199.6
$ws.Range("K91").Value = 2953.5
$ws.Range("L91").Value = 4199.6
$ws.Range("M91").Value = -1549.5
$ws.Range("N91").Value = -7007.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 756.21277
$ws.Range("I116").Value = 596.439
$ws.Range("J116").Value = 1848
$ws.Range("K116").Value = 596.439
$ws.Range("L116").Value = 1848
$ws.Range("M116").Value = 1697.561
$ws.Range("N116").Value = -6436

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1056.5588
$ws.Range("I122").Value = 788.0833
$ws.Range("J122").Value = 1700.9
$ws.Range("K122").Value = 2364.2499
$ws.Range("L122").Value = 5102.700000000001
$ws.Range("M122").Value = 85.7501000000002
$ws.Range("N122").Value = -10002.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4503.629
$ws.Range("I132").Value = 2871.302
$ws.Range("J132").Value = 14116.223
$ws.Range("K132").Value = 8613.906000000001
$ws.Range("L132").Value = 42348.669
$ws.Range("M132").Value = -6083.906000000001
$ws.Range("N132").Value = -47408.669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 756.21277
$ws.Range("I3").Value = 596.439
$ws.Range("J3").Value = 1848
$ws.Range("K3").Value = 596.439
$ws.Range("L3").Value = 1848
$ws.Range("M3").Value = -482.439
$ws.Range("N3").Value = -2076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 20000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 8682.125
$ws.Range("I82").Value = 5636.7144
$ws.Range("J82").Value = 30000
$ws.Range("K82").Value = 5636.7144
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -5253.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 8682.125
$ws.Range("I85").Value = 5636.7144
$ws.Range("J85").Value = 30000
$ws.Range("K85").Value = 5636.7144
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -4310.7144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3248.7778
$ws.Range("I86").Value = 2621.2
$ws.Range("J86").Value = 3490.1538
$ws.Range("K86").Value = 2621.2
$ws.Range("L86").Value = 3490.1538
$ws.Range("M86").Value = -1498.2
$ws.Range("N86").Value = -5736.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3248.7778
$ws.Range("I89").Value = 2621.2
$ws.Range("J89").Value = 3490.1538
$ws.Range("K89").Value = 13106
$ws.Range("L89").Value = 17450.769
$ws.Range("M89").Value = -7490
$ws.Range("N89").Value = -28682.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1960.375
$ws.Range("I107").Value = 1528.6
$ws.Range("J107").Value = 2680
$ws.Range("K107").Value = 1528.6
$ws.Range("L107").Value = 2680
$ws.Range("M107").Value = 391.4000000000001
$ws.Range("N107").Value = -6520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 5133.3335
$ws.Range("I11").Value = 2266.6667
$ws.Range("J11").Value = 8000
$ws.Range("K11").Value = 2266.6667
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = -2126.6667
$ws.Range("N11").Value = -8280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3281.7693
$ws.Range("I16").Value = 2571.5715
$ws.Range("J16").Value = 4110.3335
$ws.Range("K16").Value = 2571.5715
$ws.Range("L16").Value = 4110.3335
$ws.Range("M16").Value = -2284.5715
$ws.Range("N16").Value = -4684.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3281.7693
$ws.Range("I113").Value = 2571.5715
$ws.Range("J113").Value = 4110.3335
$ws.Range("K113").Value = 2571.5715
$ws.Range("L113").Value = 4110.3335
$ws.Range("M113").Value = -401.5715
$ws.Range("N113").Value = -8450.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3013.3333
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 3171.4285
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 9514.2855
$ws.Range("M39").Value = -2106
$ws.Range("N39").Value = -10102.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2763.1667
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2763.1667
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 8289.500100000001
$ws.Range("N62").Value = -9661.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 2763.1667
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2763.1667
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 24868.5003
$ws.Range("N65").Value = -31732.5003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3952.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3952.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 11857.5
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -13853.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 3952.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3952.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 35572.5
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -45556.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 842.5714
$ws.Range("I131").Value = 350
$ws.Range("J131").Value = 996.5
$ws.Range("K131").Value = 1050
$ws.Range("L131").Value = 2989.5
$ws.Range("M131").Value = 3990
$ws.Range("N131").Value = -13069.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 834492.7
$ws.Range("I113").Value = 1429551.8
$ws.Range("J113").Value = 1410
$ws.Range("K113").Value = 1429551.8
$ws.Range("L113").Value = 1410
$ws.Range("M113").Value = -1427381.8
$ws.Range("N113").Value = -5750

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2363.6875
$ws.Range("I122").Value = 1182.8182
$ws.Range("J122").Value = 4961.6
$ws.Range("K122").Value = 3548.4546
$ws.Range("L122").Value = 14884.8
$ws.Range("M122").Value = -1098.4546
$ws.Range("N122").Value = -19784.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 22371.258
$ws.Range("I132").Value = 47990.25
$ws.Range("J132").Value = 4287.2646
$ws.Range("K132").Value = 143970.75
$ws.Range("L132").Value = 12861.7938
$ws.Range("M132").Value = -141440.75
$ws.Range("N132").Value = -17921.7938

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2059
$ws.Range("I46").Value = 1924.5714
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 1924.5714
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -1736.5714
$ws.Range("N46").Value = -3376

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2815
$ws.Range("I68").Value = 2900
$ws.Range("J68").Value = 2802.8572
$ws.Range("K68").Value = 2900
$ws.Range("L68").Value = 2802.8572
$ws.Range("M68").Value = -2151
$ws.Range("N68").Value = -4300.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2815
$ws.Range("I71").Value = 2900
$ws.Range("J71").Value = 2802.8572
$ws.Range("K71").Value = 14500
$ws.Range("L71").Value = 14014.286
$ws.Range("M71").Value = -10756
$ws.Range("N71").Value = -21502.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 33134.8
$ws.Range("I132").Value = 25179.8
$ws.Range("J132").Value = 37112.3
$ws.Range("K132").Value = 75539.39999999999
$ws.Range("L132").Value = 111336.9
$ws.Range("M132").Value = -73009.39999999999
$ws.Range("N132").Value = -116396.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 49142.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 49142.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 49142.668
$ws.Range("N140").Value = -59502.668
